$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: column A holds a new shared string "Q9", formatted like the cells above it (A2:A10)
$ws.Cells.Item(11, 1).Value = "Q9"
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)

# Row 2
$ws.Cells.Item(2, 2).Value = 0.02846637460976744
$ws.Cells.Item(2, 3).Value = 0.9430602768606053
$ws.Cells.Item(2, 4).Value = 2.648581526649609
$ws.Cells.Item(2, 5).Value = 1.627446320666095
$ws.Cells.Item(2, 6).Value = 1.643388759852491
$ws.Cells.Item(2, 7).Value = 51

# Row 3
$ws.Cells.Item(3, 2).Value = -0.1943718813857691
$ws.Cells.Item(3, 3).Value = 0.7420322100163707
$ws.Cells.Item(3, 4).Value = 1.301007984511005
$ws.Cells.Item(3, 5).Value = 1.140617369897112
$ws.Cells.Item(3, 6).Value = 1.13534475400492
$ws.Cells.Item(3, 7).Value = 50

# Row 4
$ws.Cells.Item(4, 2).Value = -0.05603653138027285
$ws.Cells.Item(4, 3).Value = 0.7181205273575723
$ws.Cells.Item(4, 4).Value = 1.474592598453738
$ws.Cells.Item(4, 5).Value = 1.214328044003653
$ws.Cells.Item(4, 6).Value = 1.22560506124038
$ws.Cells.Item(4, 7).Value = 49

# Row 5
$ws.Cells.Item(5, 2).Value = -0.151782044080011
$ws.Cells.Item(5, 3).Value = 0.7437703319808541
$ws.Cells.Item(5, 4).Value = 1.325335475341544
$ws.Cells.Item(5, 5).Value = 1.151232155276052
$ws.Cells.Item(5, 6).Value = 1.153258924895003
$ws.Cells.Item(5, 7).Value = 48

# Row 6
$ws.Cells.Item(6, 2).Value = -0.01202353621938327
$ws.Cells.Item(6, 3).Value = 0.9685194911816654
$ws.Cells.Item(6, 4).Value = 2.135832188922473
$ws.Cells.Item(6, 5).Value = 1.461448661062876
$ws.Cells.Item(6, 6).Value = 1.47719857003534
$ws.Cells.Item(6, 7).Value = 47

# Row 7
$ws.Cells.Item(7, 2).Value = -0.1554833155668876
$ws.Cells.Item(7, 3).Value = 0.8002215085028028
$ws.Cells.Item(7, 4).Value = 1.375155290479728
$ws.Cells.Item(7, 5).Value = 1.172670154169419
$ws.Cells.Item(7, 6).Value = 1.175160419656847
$ws.Cells.Item(7, 7).Value = 46

# Row 8
$ws.Cells.Item(8, 2).Value = 0.04659670587843785
$ws.Cells.Item(8, 3).Value = 0.8050767692874348
$ws.Cells.Item(8, 4).Value = 1.154340410675184
$ws.Cells.Item(8, 5).Value = 1.074402350460564
$ws.Cells.Item(8, 6).Value = 1.085520529677319
$ws.Cells.Item(8, 7).Value = 45

# Row 9
$ws.Cells.Item(9, 2).Value = -0.1084043684214821
$ws.Cells.Item(9, 3).Value = 0.735808405588202
$ws.Cells.Item(9, 4).Value = 1.20283139358349
$ws.Cells.Item(9, 5).Value = 1.096736702031755
$ws.Cells.Item(9, 6).Value = 1.1039834322736
$ws.Cells.Item(9, 7).Value = 44

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1182394667038753
$ws.Cells.Item(10, 3).Value = 1.047860666370489
$ws.Cells.Item(10, 4).Value = 2.217418452168737
$ws.Cells.Item(10, 5).Value = 1.489099879849816
$ws.Cells.Item(10, 6).Value = 1.501965607916916
$ws.Cells.Item(10, 7).Value = 43

# Row 11
$ws.Cells.Item(11, 2).Value = -0.08993503749128429
$ws.Cells.Item(11, 3).Value = 0.8551623832460461
$ws.Cells.Item(11, 4).Value = 1.503988522260922
$ws.Cells.Item(11, 5).Value = 1.226372097799409
$ws.Cells.Item(11, 6).Value = 1.237895626577417
$ws.Cells.Item(11, 7).Value = 42
